# Daily attendance processing - normalize "Recorded By" entries so that the
# literal "System" token is reordered to sit after the other recorder(s)
# it was paired with, matching the upstream export format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# Locate the "Recorded By" column dynamically from the header row,
# falling back to column 7 (G) if it can't be found.
$recordedByCol = 0
for ($c = 0; $c -lt $colCount; $c++) {
    $header = $ws.Cells.Item($firstRow, $firstCol + $c).Text
    if ($header -eq "Recorded By") {
        $recordedByCol = $firstCol + $c
        break
    }
}
if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

for ($i = 1; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
